# Fruta / hortaliza, semanal
# Insert a new weekly record at row 43 (pushing the existing rows 43-57 down to 44-58)
# and populate it with the new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 43:57 down by one row to make room for the new record.
$ws.Rows.Item(43).Insert()

# Fill in the new row with the new weekly data point.
$ws.Cells.Item(43, 1).Value2 = 6
$ws.Cells.Item(43, 2).Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(43, 3).Value2 = "Metropolitana"
$ws.Cells.Item(43, 4).Value2 = 44559
$ws.Cells.Item(43, 5).Value2 = 13
$ws.Cells.Item(43, 6).Value2 = "Fruta"
$ws.Cells.Item(43, 7).Value2 = 100101
$ws.Cells.Item(43, 8).Value2 = "Berries"
$ws.Cells.Item(43, 9).Value2 = 100101008
$ws.Cells.Item(43, 10).Value2 = "Mora"
$ws.Cells.Item(43, 11).Value2 = "Sin especificar"
$ws.Cells.Item(43, 12).Value2 = "Primera"
$ws.Cells.Item(43, 13).Value2 = 200
$ws.Cells.Item(43, 14).Value2 = 6000
$ws.Cells.Item(43, 15).Value2 = 6000
$ws.Cells.Item(43, 16).Value2 = 6000
$ws.Cells.Item(43, 17).Value2 = "`$/bandeja 2 kilos"
$ws.Cells.Item(43, 18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(43, 19).Value2 = 3000
$ws.Cells.Item(43, 20).Value2 = 2
